# "Generate Report for Handoff"
# The localization-status report is refreshed: the file
# 755a219d-5d8b-44d5-8154-e574f1dc4ff8.md (row 3 on each sheet) has now
# been handed off for translation, so its status flips from
# "In Translation" to "Ready for handoff" on the Overview sheet and on
# each per-locale sheet, and the per-locale "Latest Handoff Datetime" is
# updated to the new handoff timestamp.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: update status for both locales on the 755a219d... row
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# zh-cn sheet: update status and handoff datetime
$wsZhCn.Range("B3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "2016-03-10 07:59:33"

# de-de sheet: update status and handoff datetime
$wsDeDe.Range("B3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "2016-03-10 07:59:42"
